$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 37, pushing the existing rows 37-80
# down to 38-81 (this also extends the used range / dimension to A1:R81).
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new weekly price record.
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44539
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112052
$ws.Range("G37").Value = "Albahaca"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 60
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 8000
$ws.Range("M37").Value = 8000
$ws.Range("N37").Value = "$/docena de matas"
$ws.Range("O37").Value = "Región Metropolitana"
$ws.Range("P37").Value = 1333
$ws.Range("Q37").Value = 6
$ws.Range("R37").Value = "Hortaliza"

# Make sure the date cell keeps the same custom date number format used by
# every other row in column D.
$ws.Range("D37").NumberFormat = $ws.Range("D38").NumberFormat
